# The "Artfynd" sheet rows 2-16 each describe one species observation. This
# edit reshuffles which observation sits in which row (a pure row-data
# permutation - same 15 records, new row order), matching the upstream
# source export order.
#
# Only the columns whose value actually differs between the old and the new
# occupant of a row are written, and only cells whose value truly changes
# are touched at all - columns that are identical across every one of the 15
# rows (C,P,S,T,U,V,W,Z,AB,AD,AE,AG,AT,AW,AY) are left completely alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row permutation: new row R ends up holding the data that currently lives
# in row Map[R].
$map = @{}
$map[2]  = 6
$map[3]  = 7
$map[4]  = 8
$map[5]  = 9
$map[6]  = 2
$map[7]  = 10
$map[8]  = 11
$map[9]  = 12
$map[10] = 13
$map[11] = 14
$map[12] = 15
$map[13] = 16
$map[14] = 3
$map[15] = 4
$map[16] = 5

# Columns whose value can differ from row to row.
$numCols  = @("A","B","E","Q","R")
$textCols = @("D","F","G","H","I","Y","AA","AC","AX")
$allCols  = $numCols + $textCols

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16)

# Snapshot every current value BEFORE any write, so the permutation can be
# computed purely from memory (a source row's original data is never read
# back after being overwritten).
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $allCols) { $rowData[$c] = $ws.Range("$c$r").Value() }
    $snapshot[$r] = $rowData
}

# Write the permuted values back - but only where the value actually
# changes, and only flip a cell to text formatting when it is about to hold
# a numeric-/date-looking string (column I holds plain digits like "3"; Y
# and AA hold "yyyy-mm-dd" strings) so Excel's Value setter doesn't silently
# re-type it as a Double/Date.
$forceTextCols = @{ "I" = $true; "Y" = $true; "AA" = $true }

foreach ($r in $rows) {
    $src = $map[$r]
    if ($src -eq $r) { continue }
    $data = $snapshot[$src]
    $cur = $snapshot[$r]
    foreach ($c in $numCols) {
        if ($data[$c] -ne $cur[$c]) {
            $ws.Range("$c$r").Value() = $data[$c]
        }
    }
    foreach ($c in $textCols) {
        if ([string]$data[$c] -ne [string]$cur[$c]) {
            if ($forceTextCols.ContainsKey($c)) {
                $ws.Range("$c$r").NumberFormat = "@"
            }
            $ws.Range("$c$r").Value() = [string]$data[$c]
        }
    }
}
